$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "Implied weights are same with Black Litterman weights since no view has been supplied yet"
$ws.Range("B15").Select()
